# ---------------------------------------------------------------------------
# B6-PowerPoint.pptx edit
#
# 1) The deck's theme colour scheme (ppt/theme/theme1.xml, used by the main
#    slide master -> every slide) is switched from the "Integral" / "Red
#    Violet" palette to the standard "Office" palette (the palette that used
#    to live only in ppt/theme/theme2.xml, which is used by the notes
#    master). PowerPoint's object model exposes this through
#    Master.ColorScheme.Colors(i).RGB (the classic 12-slot scheme: dk1, lt1,
#    dk2, lt2, accent1-6, hlink, folHlink), so we rewrite every slot to the
#    "Office" theme's RGB values.
#
# 2) Three tables (one each on slide 14, slide 15 and slide 16) get their
#    table style switched from {22A4C0EB-7B67-4EE0-89AE-B1121312A1B2} to
#    {8AAD3AF0-AF0F-4EEE-A877-3086AC66DAD9} via Table.ApplyStyle(guid).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Recolour the master theme to the "Office" palette -------------------
# RGBColor.RGB uses the classic COM 0x00BBGGRR packing, i.e. R + G*256 + B*65536.
$officeColors = @{
    1  = 0          # dk1      000000
    2  = 16777215   # lt1      FFFFFF
    3  = 6968388    # dk2      44546A
    4  = 15132391   # lt2      E7E6E6
    5  = 13998939   # accent1  5B9BD5
    6  = 3243501    # accent2  ED7D31
    7  = 10855845   # accent3  A5A5A5
    8  = 49407      # accent4  FFC000
    9  = 12874308   # accent5  4472C4
    10 = 4697456    # accent6  70AD47
    11 = 12673797   # hlink    0563C1
    12 = 7491477    # folHlink 954F72
}

$masterScheme = $p.SlideMaster.ColorScheme
foreach ($slot in $officeColors.Keys) {
    $masterScheme.Colors($slot).RGB = $officeColors[$slot]
}

# --- 2) Re-style the three tables on slides 14-16 ---------------------------
$newStyleId = "{8AAD3AF0-AF0F-4EEE-A877-3086AC66DAD9}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
